$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LoginCredentials")
$ws.Range("A1").Font.Bold = $true
$ws.Range("A1").Font.Size = 12
$ws.Range("A1").Font.Name = "Helvetica Neue"
$ws.Range("A1").Font.Color = 13408767
